$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs logically
# before the current row 112, so insert a blank row there which pushes the
# existing rows 112-172 down to 113-173 (dimension grows to A1:T173).
$ws.Rows.Item(112).Insert()

# Fill in the newly inserted row 112 with the new record's data.
$ws.Range("A112").Value = 4
$ws.Range("B112").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C112").Value = "Los Lagos"
$ws.Range("D112").Value = 44523
$ws.Range("E112").Value = 10
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100103
$ws.Range("H112").Value = "Frutos de hueso (carozo)"
$ws.Range("I112").Value = 100103004
$ws.Range("J112").Value = "Durazno"
$ws.Range("K112").Value = "Florida King"
$ws.Range("L112").Value = "Segunda"
$ws.Range("M112").Value = 600
$ws.Range("N112").Value = 13000
$ws.Range("O112").Value = 13500
$ws.Range("P112").Value = 13250
$ws.Range("Q112").Value = "$/caja 14 kilos empedrada"
$ws.Range("R112").Value = "Provincia de Limarí"
$ws.Range("S112").Value = 946
$ws.Range("T112").Value = 14
